$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every touched cell to Text format first so numeric-looking
# strings (e.g. "173.60", "0.0000246") round-trip verbatim instead of
# being coerced to a Number by COM Value auto-detection (which would
# drop trailing zeros / switch to scientific notation).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.959.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.32%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.110.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.49%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.56%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.65%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.104.84"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.47%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.67%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.45%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.478"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.11"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.99%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.630.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.984.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.116.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.61%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "477.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.72"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.79"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.87%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.92%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.94"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.68%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.96"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.58%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.59"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0965"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.74%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.61%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.977"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.46"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.08"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.48%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.04"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.309"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.98%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.59"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.795.58"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0354"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "378.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.42%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -11.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.80%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.70"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.14%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.76%  "
